# "some questions by 9/10"
# Adds a new "Sheet2" with a numbered question list (col A = row number,
# col B = question text for the first few rows), adjusts Sheet1's
# selection/formatting a bit, and makes Sheet2 the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 0. Small formatting touch-ups on Sheet1 (do this before Sheet2 becomes
#    active so Sheet2 stays the final active/visible tab)
# ---------------------------------------------------------------------
$ws1.Range("A1:A6").RowHeight = 14.25
$ws1.Range("A1").ColumnWidth = 94.125
[void]$ws1.Range("A6").Select()

# ---------------------------------------------------------------------
# 1. Create Sheet2 right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# 2. Fill Sheet2 data
#    Column A: row numbers 1..19
#    Column B: the first five rows hold the new questions
# ---------------------------------------------------------------------
for ($i = 1; $i -le 19; $i++) {
    $ws2.Cells.Item($i, 1).Value = $i
}

$ws2.Range("B1").Value = "pixel spacing对分析的影响"
$ws2.Range("B2").Value = "一般图像分析的方法:CNN,DNN?"
$ws2.Range("B3").Value = "图像需要做什么样的data preprocessing吗? 例如pixel normalization之类的?"
$ws2.Range("B4").Value = "要不要先做分类分析 然后再具体分析+框"
$ws2.Range("B5").Value = "activation function等之类的需要什么注意的?"

# Match the look & feel (font/alignment) used for the question column on
# Sheet1 by copying its formatting over to B1 only (B2:B5 stay unstyled).
$ws1.Range("A1").Copy()
[void]$ws2.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The first question mixes two fonts in a single cell: "pixel spacing" in
# the regular font, and the Chinese remainder in 宋体 (SimSun).
$chars = $ws2.Range("B1").Characters(14, 6)
$chars.Font.Name = "宋体"

# ---------------------------------------------------------------------
# 3. Sheet2 view: select B6, make Sheet2 the active/visible tab
# ---------------------------------------------------------------------
$ws2.Activate()
[void]$ws2.Range("B6").Select()
